$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 16675810
$ws.Range("I69").Value = 33340654
$ws.Range("J69").Value = 10966.667
$ws.Range("K69").Value = 100021962
$ws.Range("L69").Value = 32900.001
$ws.Range("M69").Value = -100021088
$ws.Range("N69").Value = -34648.001

$ws.Range("H70").Value = 813
$ws.Range("I70").Value = 623.25
$ws.Range("J70").Value = 1066
$ws.Range("K70").Value = 1869.75
$ws.Range("L70").Value = 3198
$ws.Range("M70").Value = -1599.75
$ws.Range("N70").Value = -3738

$ws.Range("H72").Value = 16675810
$ws.Range("I72").Value = 33340654
$ws.Range("J72").Value = 10966.667
$ws.Range("K72").Value = 300065886
$ws.Range("L72").Value = 98700.003
$ws.Range("M72").Value = -300061518
$ws.Range("N72").Value = -107436.003

$ws.Range("H73").Value = 813
$ws.Range("I73").Value = 623.25
$ws.Range("J73").Value = 1066
$ws.Range("K73").Value = 1869.75
$ws.Range("L73").Value = 3198
$ws.Range("M73").Value = -933.75
$ws.Range("N73").Value = -5070

$ws.Range("H80").Value = 16667120
$ws.Range("I80").Value = 307
$ws.Range("J80").Value = 25000526
$ws.Range("K80").Value = 921
$ws.Range("L80").Value = 75001578
$ws.Range("M80").Value = 77
$ws.Range("N80").Value = -75003574

$ws.Range("H83").Value = 16667120
$ws.Range("I83").Value = 307
$ws.Range("J83").Value = 25000526
$ws.Range("K83").Value = 2763
$ws.Range("L83").Value = 225004734
$ws.Range("M83").Value = 2229
$ws.Range("N83").Value = -225014718

$ws.Range("H132").Value = 2141.1428
$ws.Range("I132").Value = 2249.75
$ws.Range("J132").Value = 1996.3334
$ws.Range("K132").Value = 6749.25
$ws.Range("L132").Value = 5989.0002
$ws.Range("M132").Value = -4219.25
$ws.Range("N132").Value = -11049.0002

$ws.Range("H137").Value = 1288469.5
$ws.Range("I137").Value = 6404.4062
$ws.Range("J137").Value = 2180341
$ws.Range("K137").Value = 19213.2186
$ws.Range("L137").Value = 6541023
$ws.Range("M137").Value = -16663.2186
$ws.Range("N137").Value = -6546123

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 2653.3333
$ws.Range("I16").Value = 2653.3333
$ws.Range("K16").Value = 2653.3333
$ws.Range("M16").Value = -2366.3333

$ws.Range("H61").Value = 1092558.4
$ws.Range("I61").Value = 30422.475
$ws.Range("J61").Value = 3466744.5
$ws.Range("K61").Value = 30422.475
$ws.Range("L61").Value = 3466744.5
$ws.Range("M61").Value = -30210.475
$ws.Range("N61").Value = -3467168.5

$ws.Range("H98").Value = 77500
$ws.Range("J98").Value = 77500
$ws.Range("L98").Value = 77500
$ws.Range("N98").Value = -83490

$ws.Range("H103").Value = 69463.664
$ws.Range("J103").Value = 69463.664
$ws.Range("L103").Value = 69463.664
$ws.Range("N103").Value = -71807.664

$ws.Range("H122").Value = 1861.04
$ws.Range("I122").Value = 1296.7727
$ws.Range("J122").Value = 5999
$ws.Range("K122").Value = 3890.3181
$ws.Range("L122").Value = 17997
$ws.Range("M122").Value = -1440.3181
$ws.Range("N122").Value = -22897

$ws.Range("H136").Value = 1092558.4
$ws.Range("I136").Value = 30422.475
$ws.Range("J136").Value = 3466744.5
$ws.Range("K136").Value = 91267.42499999999
$ws.Range("L136").Value = 10400233.5
$ws.Range("M136").Value = -88717.42499999999
$ws.Range("N136").Value = -10405333.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3885.6875
$ws.Range("I86").Value = 2855.25
$ws.Range("J86").Value = 5603.0835
$ws.Range("K86").Value = 2855.25
$ws.Range("L86").Value = 5603.0835
$ws.Range("M86").Value = -1732.25
$ws.Range("N86").Value = -7849.0835

$ws.Range("H89").Value = 3885.6875
$ws.Range("I89").Value = 2855.25
$ws.Range("J89").Value = 5603.0835
$ws.Range("K89").Value = 14276.25
$ws.Range("L89").Value = 28015.4175
$ws.Range("M89").Value = -8660.25
$ws.Range("N89").Value = -39247.4175

$ws.Range("H99").Value = 7219
$ws.Range("I99").Value = 8871
$ws.Range("J99").Value = 2971
$ws.Range("K99").Value = 8871
$ws.Range("L99").Value = 2971
$ws.Range("M99").Value = -7373
$ws.Range("N99").Value = -5967

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 14299591
$ws.Range("I16").Value = 20410844
$ws.Range("J16").Value = 39999.332
$ws.Range("K16").Value = 20410844
$ws.Range("L16").Value = 39999.332
$ws.Range("M16").Value = -20410557
$ws.Range("N16").Value = -40573.332

$ws.Range("H31").Value = 4959.875
$ws.Range("J31").Value = 5142.3066
$ws.Range("L31").Value = 5142.3066
$ws.Range("N31").Value = -5732.3066

$ws.Range("H34").Value = 4959.875
$ws.Range("J34").Value = 5142.3066
$ws.Range("L34").Value = 5142.3066
$ws.Range("N34").Value = -5546.3066

$ws.Range("H94").Value = 1462.2
$ws.Range("I94").Value = 1032.3334
$ws.Range("J94").Value = 1748.7778
$ws.Range("K94").Value = 1032.3334
$ws.Range("L94").Value = 1748.7778
$ws.Range("M94").Value = -581.3334
$ws.Range("N94").Value = -2650.7778

$ws.Range("H105").Value = 3156.5386
$ws.Range("I105").Value = 2502.7778
$ws.Range("J105").Value = 4627.5
$ws.Range("K105").Value = 2502.7778
$ws.Range("L105").Value = 4627.5
$ws.Range("M105").Value = -755.7777999999998
$ws.Range("N105").Value = -8121.5

$ws.Range("H113").Value = 14299591
$ws.Range("I113").Value = 20410844
$ws.Range("J113").Value = 39999.332
$ws.Range("K113").Value = 20410844
$ws.Range("L113").Value = 39999.332
$ws.Range("M113").Value = -20408674
$ws.Range("N113").Value = -44339.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1251.4814
$ws.Range("I5").Value = 811.38464
$ws.Range("J5").Value = 1660.1428
$ws.Range("K5").Value = 2434.15392
$ws.Range("L5").Value = 4980.428400000001
$ws.Range("M5").Value = -2322.15392
$ws.Range("N5").Value = -5204.428400000001

$ws.Range("H11").Value = 602.2632
$ws.Range("I11").Value = 655.4
$ws.Range("J11").Value = 403
$ws.Range("K11").Value = 1966.2
$ws.Range("L11").Value = 1209
$ws.Range("M11").Value = -1826.2
$ws.Range("N11").Value = -1489

$ws.Range("H15").Value = 618.9231
$ws.Range("I15").Value = 635.5454999999999
$ws.Range("J15").Value = 527.5
$ws.Range("K15").Value = 1906.6365
$ws.Range("L15").Value = 1582.5
$ws.Range("M15").Value = -1766.6365
$ws.Range("N15").Value = -1862.5

$ws.Range("H68").Value = 2780.3157
$ws.Range("I68").Value = 2258.3333
$ws.Range("J68").Value = 3675.1428
$ws.Range("K68").Value = 6774.999899999999
$ws.Range("L68").Value = 11025.4284
$ws.Range("M68").Value = -5963.999899999999
$ws.Range("N68").Value = -12647.4284

$ws.Range("H71").Value = 2780.3157
$ws.Range("I71").Value = 2258.3333
$ws.Range("J71").Value = 3675.1428
$ws.Range("K71").Value = 20324.9997
$ws.Range("L71").Value = 33076.2852
$ws.Range("M71").Value = -16268.9997
$ws.Range("N71").Value = -41188.2852

$ws.Range("H135").Value = 1251.4814
$ws.Range("I135").Value = 811.38464
$ws.Range("J135").Value = 1660.1428
$ws.Range("K135").Value = 7302.46176
$ws.Range("L135").Value = 14941.2852
$ws.Range("M135").Value = -4767.46176
$ws.Range("N135").Value = -20011.2852

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 125000
$ws.Range("J93").Value = 125000
$ws.Range("L93").Value = 125000
$ws.Range("N93").Value = -128744

$ws.Range("H126").Value = 3396.5217
$ws.Range("I126").Value = 3339.4443
$ws.Range("J126").Value = 3602
$ws.Range("K126").Value = 10018.3329
$ws.Range("L126").Value = 10806
$ws.Range("M126").Value = -7548.332900000001
$ws.Range("N126").Value = -15746

$ws.Range("H135").Value = 119999.25
$ws.Range("J135").Value = 119999.25
$ws.Range("L135").Value = 119999.25
$ws.Range("N135").Value = -130139.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 204699.8
$ws.Range("I81").Value = 6750
$ws.Range("J81").Value = 336666.34
$ws.Range("K81").Value = 13500
$ws.Range("L81").Value = 673332.6800000001
$ws.Range("M81").Value = -12439
$ws.Range("N81").Value = -675454.6800000001

$ws.Range("H84").Value = 204699.8
$ws.Range("I84").Value = 6750
$ws.Range("J84").Value = 336666.34
$ws.Range("K84").Value = 67500
$ws.Range("L84").Value = 3366663.4
$ws.Range("M84").Value = -62196
$ws.Range("N84").Value = -3377271.4

$ws.Range("H122").Value = 4423.4375
$ws.Range("I122").Value = 3806.9092
$ws.Range("J122").Value = 5779.8
$ws.Range("K122").Value = 11420.7276
$ws.Range("L122").Value = 17339.4
$ws.Range("M122").Value = -8970.7276
$ws.Range("N122").Value = -22239.4
